# Commit: "Cập nhật lại kết quả test với code DT ANN mới"
# (Update test results with the new DT ANN code)
#
# The workbook tab "Sheet2" (internally xl/worksheets/sheet6.xml) holds two
# blocks of per-stock accuracy figures (rows 2-12 and rows 15-25), each
# followed by an AVERAGE summary row (13 and 26). The edit refreshes all of
# the raw figures with newly produced numbers and turns the summary cells
# from hard-coded numbers into live =AVERAGE(...) formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# ---- Block 1: rows 2-12 (years 2003 / 2004 / 2005 in columns B / D / F) ----
$ws.Range("B2").Value  = 83.5
$ws.Range("D2").Value  = 78.5
$ws.Range("F2").Value  = 80

$ws.Range("B3").Value  = 74
$ws.Range("D3").Value  = 74.5
$ws.Range("F3").Value  = 72

$ws.Range("B4").Value  = 85
$ws.Range("D4").Value  = 64.5
$ws.Range("F4").Value  = 60

$ws.Range("B5").Value  = 83.5
$ws.Range("D5").Value  = 84.5
$ws.Range("F5").Value  = 85.5

$ws.Range("B6").Value  = 92
$ws.Range("D6").Value  = 89
$ws.Range("F6").Value  = 89

$ws.Range("B7").Value  = 82
$ws.Range("D7").Value  = 66.5
$ws.Range("F7").Value  = 63

$ws.Range("B8").Value  = 95
$ws.Range("D8").Value  = 95
$ws.Range("F8").Value  = 93.5

$ws.Range("B9").Value  = 83
$ws.Range("D9").Value  = 87
$ws.Range("F9").Value  = 81.5

$ws.Range("B10").Value = 83
$ws.Range("D10").Value = 87
$ws.Range("F10").Value = 84.5

$ws.Range("B11").Value = 88
$ws.Range("D11").Value = 87.5
$ws.Range("F11").Value = 85.5

$ws.Range("B12").Value = 91
$ws.Range("D12").Value = 92
$ws.Range("F12").Value = 92

# Row 13 used to hold static averages -> make them real formulas.
$ws.Range("B13").Formula = "=AVERAGE(B2:B12)"
$ws.Range("D13").Formula = "=AVERAGE(D2:D12)"
$ws.Range("F13").Formula = "=AVERAGE(F2:F12)"

# ---- Block 2: rows 15-25 ----
$ws.Range("B15").Value = 82.5
$ws.Range("D15").Value = 80
$ws.Range("F15").Value = 82.5

$ws.Range("B16").Value = 64
$ws.Range("D16").Value = 62.5
$ws.Range("F16").Value = 64

$ws.Range("B17").Value = 80
$ws.Range("D17").Value = 80
$ws.Range("F17").Value = 60

$ws.Range("B18").Value = 83
$ws.Range("D18").Value = 81
$ws.Range("F18").Value = 84

$ws.Range("B19").Value = 84
$ws.Range("D19").Value = 78.5
$ws.Range("F19").Value = 78.5

$ws.Range("B20").Value = 70.5
$ws.Range("D20").Value = 54.5
$ws.Range("F20").Value = 54

$ws.Range("B21").Value = 90
$ws.Range("D21").Value = 91.5
$ws.Range("F21").Value = 92

$ws.Range("B22").Value = 82
$ws.Range("D22").Value = 81.5
$ws.Range("F22").Value = 70.5

$ws.Range("B23").Value = 82
$ws.Range("D23").Value = 73
$ws.Range("F23").Value = 67

$ws.Range("B24").Value = 86.5
$ws.Range("D24").Value = 81.5
$ws.Range("F24").Value = 82

$ws.Range("B25").Value = 84
$ws.Range("D25").Value = 83.5
$ws.Range("F25").Value = 85

# Row 26 used to hold static averages -> make them real formulas.
$ws.Range("B26").Formula = "=AVERAGE(B15:B25)"
$ws.Range("D26").Formula = "=AVERAGE(D15:D25)"
$ws.Range("F26").Formula = "=AVERAGE(F15:F25)"

# Column A got a touch wider once the longer labels were re-entered.
$ws.Columns.Item(1).ColumnWidth = 11.83

# Scroll/selection position left by the author after the last edit.
$ws.Range("F26").Select()
